$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) sometimes holds plain-looking decimal strings
# (e.g. "1.010", "8.120"). Left alone, Excel auto-converts those to
# numbers on assignment and drops the significant trailing zero, so we
# force the cell to Text format first for just those values, which
# keeps the literal string (matches the source data which is all text).

$ws.Range("D2").Value = '28.022.56'
$ws.Range("E2").Value = '  +0.24%  '

$ws.Range("D3").Value = '1.846.32'
$ws.Range("E3").Value = '  -0.88%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.011'
$ws.Range("E4").Value = '  +0.86%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '331.77'
$ws.Range("E5").Value = '  -1.30%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.010'
$ws.Range("E6").Value = '  +0.71%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4519'
$ws.Range("E7").Value = '  -3.90%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3894'
$ws.Range("E8").Value = '  +0.21%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.56'
$ws.Range("E9").Value = '  +1.48%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07739'
$ws.Range("E10").Value = '  -2.95%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9756'
$ws.Range("E11").Value = '  -0.14%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.12'
$ws.Range("E12").Value = '  -1.53%  '

$ws.Range("D13").Value = '1.842.33'
$ws.Range("E13").Value = '  -0.56%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.767'
$ws.Range("E14").Value = '  -2.69%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.973'
$ws.Range("E15").Value = '  -3.19%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.011'
$ws.Range("E16").Value = '  +0.78%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '86.87'
$ws.Range("E17").Value = '  -5.19%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06533'
$ws.Range("E18").Value = '  -1.41%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.00001015'
$ws.Range("E19").Value = '  -2.09%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.88'
$ws.Range("E20").Value = '  -3.20%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.008'
$ws.Range("E21").Value = '  +0.50%  '

$ws.Range("D22").Value = '28.041.61'
$ws.Range("E22").Value = '  +0.36%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.297'
$ws.Range("E23").Value = '  -1.66%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.58'
$ws.Range("E24").Value = '  -2.93%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.259'
$ws.Range("E25").Value = '  -1.49%  '

$ws.Range("D26").Value = '2.075.66'
$ws.Range("E26").Value = '  +0.43%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '154.92'
$ws.Range("E27").Value = '  -2.51%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.09'
$ws.Range("E28").Value = '  -2.37%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.027'
$ws.Range("E29").Value = '  -3.19%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.213'
$ws.Range("E30").Value = '  -4.27%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '116.05'
$ws.Range("E31").Value = '  -2.73%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09240'
$ws.Range("E32").Value = '  -2.50%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9292'
$ws.Range("E33").Value = '  -2.73%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.615'
$ws.Range("E34").Value = '  +0.98%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.364'
$ws.Range("E35").Value = '  +1.59%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.167'
$ws.Range("E36").Value = '  -2.58%  '

$ws.Range("E37").Value = '  -1.57%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02176'
$ws.Range("E38").Value = '  -2.53%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.120'
$ws.Range("E39").Value = '  -2.02%  '

$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.143'
$ws.Range("E40").Value = '  -1.60%  '

$ws.Range("B41").Value = 'Frax'
$ws.Range("C41").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.009'
$ws.Range("E41").Value = '  +0.70%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5641'
$ws.Range("E42").Value = '  -4.26%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1781'
$ws.Range("E43").Value = '  -4.32%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '9.875'
$ws.Range("E44").Value = '  -2.90%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.238'
$ws.Range("E45").Value = '  -3.46%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.258'
$ws.Range("E46").Value = '  +23.99%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.07168'

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5351'
$ws.Range("E48").Value = '  -3.01%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '11.71'
$ws.Range("E49").Value = '  -3.30%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.858'
$ws.Range("E50").Value = '  -4.60%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '109.43'
$ws.Range("E51").Value = '  -1.89%  '
